$d = $word.ActiveDocument

# The document generator split these two text runs without marking
# them as space-preserving, even though their text starts/ends with
# significant whitespace (a non-breaking space before the colon, and a
# trailing space). Re-applying the exact same text to each run makes
# Word mark the <w:t> element with xml:space="preserve", which is the
# fix described by the commit ("Fixed space preserve after parser
# spliting.").
#
# Both strings use a non-breaking space (U+00A0) before the colon,
# matching the original document content exactly.
$nbsp = [char]0x00A0

$tBookmark   = "Test bookmark" + $nbsp + ": "
$tLinkBefore = "Test link before bookmark" + $nbsp + ": "

$d.Content.Find.Execute($tBookmark, $true, $false, $false, $false, $false,
                         $true, 1, $false, $tBookmark, 2)

$d.Content.Find.Execute($tLinkBefore, $true, $false, $false, $false, $false,
                         $true, 1, $false, $tLinkBefore, 2)
